$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F5").Value = 1
$ws.Range("F7").Value = 2
$ws.Range("F10").Value = 2
$ws.Range("F13").Value = -1
$ws.Range("F16").Value = -6
$ws.Range("F19").Value = 7
$ws.Range("F22").Value = 8
$ws.Range("F23").Value = -1
$ws.Range("F26").Value = -9
$ws.Range("F27").Value = -5
$ws.Range("F31").Value = 4

$wb.Save()
